$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 194, shifting existing rows 194:223 down to 195:224
$ws.Rows.Item(194).Insert()

# Populate the new row 194 with the new weekly price record
$ws.Cells.Item(194, 1).Value = 8
$ws.Cells.Item(194, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(194, 3).Value = "Coquimbo"
$ws.Cells.Item(194, 4).Value = "2023-04-18"
$ws.Cells.Item(194, 5).Value = 4
$ws.Cells.Item(194, 6).Value = 100112040
$ws.Cells.Item(194, 7).Value = "Cilantro"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 2400
$ws.Cells.Item(194, 11).Value = 2000
$ws.Cells.Item(194, 12).Value = 2500
$ws.Cells.Item(194, 13).Value = 2250
$ws.Cells.Item(194, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(194, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(194, 16).Value = 1500
$ws.Cells.Item(194, 17).Value = 1.5
$ws.Cells.Item(194, 18).Value = "Hortaliza"
